$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.894.32"
$ws.Range("E2").Value = "  +2.98%  "

$ws.Range("D3").Value = "1.883.25"
$ws.Range("E3").Value = "  +2.98%  "

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "324.34"
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").Value = "0.4676"

$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("D9").Value = "0.07920"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").Value = "0.9836"
$ws.Range("E10").Value = "  +2.35%  "

$ws.Range("D11").Value = "22.38"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").Value = "1.922.93"
$ws.Range("E12").Value = "  +5.20%  "

$ws.Range("D13").Value = "5.752"
$ws.Range("E13").Value = "  +1.61%  "

$ws.Range("D14").Value = "7.017"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "0.06983"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").Value = "88.87"
$ws.Range("E16").Value = "  +2.63%  "

$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "0.00001011"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "16.99"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "28.892.51"
$ws.Range("E21").Value = "  +2.85%  "

$ws.Range("D22").Value = "5.355"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("D23").Value = "11.10"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("D24").Value = "2.121"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("D25").Value = "2.093.15"
$ws.Range("E25").Value = "  +2.32%  "

$ws.Range("D26").Value = "153.51"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("D27").Value = "19.38"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").Value = "5.802"
$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("E29").Value = "  +1.52%  "

$ws.Range("D30").Value = "119.84"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("D32").Value = "0.9414"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("E34").Value = "  +3.18%  "

$ws.Range("D35").Value = "3.346"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").Value = "0.05932"
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").Value = "0.02128"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").Value = "1.162"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").Value = "7.895"
$ws.Range("E39").Value = "  +3.22%  "

$ws.Range("E40").Value = "  +2.54%  "

$ws.Range("D41").Value = "0.1801"
$ws.Range("E41").Value = "  +1.87%  "

$ws.Range("D42").Value = "10.01"
$ws.Range("E42").Value = "  +0.85%  "

$ws.Range("D43").Value = "0.07315"
$ws.Range("E43").Value = "  +4.44%  "

$ws.Range("D44").Value = "11.82"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.176"
$ws.Range("E45").Value = "  -2.59%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5357"
$ws.Range("E46").Value = "  +1.59%  "

$ws.Range("D47").Value = "2.127"
$ws.Range("E47").Value = "  -4.78%  "

$ws.Range("D48").Value = "1.850"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D49").Value = "114.23"
$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D50").Value = "2.372"
$ws.Range("E50").Value = "  +2.63%  "

$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.54%  "
